# Sets the birthday (column F) values for the MOCK_DATA3 table on "Feuil1",
# replacing the placeholder "recently generated" dates with real birthdate
# serials, and applies a date number format to the F1 header cell (selecting
# the whole column first, the way Excel does when you click the column
# header and then apply Format > Cells > Date).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New birthday values (date serials) for rows 2-82 (column F).
$birthdays = @(30946, 28174, 33321, 36799, 27310, 29066, 23795, 24610, 30272, 33421, 29859, 32783, 22754, 31948, 22511, 29714, 34144, 34645, 34878, 29651, 32605, 33052, 24946, 26241, 33806, 28062, 36467, 24094, 27880, 33952, 25285, 26764, 32719, 29603, 33144, 22708, 31980, 24602, 28301, 35956, 32702, 36408, 36328, 22809, 24082, 32495, 36229, 25795, 33365, 34805, 34052, 28898, 21977, 33726, 29074, 26031, 32358, 35342, 30440, 32663, 31508, 31865, 35327, 23303, 24260, 28242, 25008, 23902, 28051, 23540, 36122, 31632, 27597, 22922, 25616, 34417, 25646, 22323, 24021, 34861, 25772)

for ($i = 0; $i -lt $birthdays.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $birthdays[$i]
}

# Select column F (mirrors clicking the column header) and apply a date
# number format to it; this mints a new style for the F1 header cell
# (same green header font/fill, now with the date number format too).
$ws.Columns("F:F").Select()
$ws.Range("F1").NumberFormat = "mm-dd-yy"
